$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3262
$ws1.Range("F3").Value = 738
$ws1.Range("F5").Value = 6898
$ws1.Range("F6").Value = 2159
$ws1.Range("F8").Value = 81
$ws1.Range("F12").Value = 26
$ws1.Range("F14").Value = 194
$ws1.Range("F15").Value = 39

# Sheet "全部类型" (fourth sheet) - same updates, rows shifted by +1 due to an
# extra entry present only in this consolidated sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3262
$ws4.Range("F4").Value = 738
$ws4.Range("F6").Value = 6898
$ws4.Range("F7").Value = 2159
$ws4.Range("F9").Value = 81
$ws4.Range("F13").Value = 26
$ws4.Range("F15").Value = 194
$ws4.Range("F16").Value = 39
